# Add a new "Expression" worksheet after the existing "Users" sheet,
# to save expressions for the session.
$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")
$ws = $wb.Worksheets.Add($null, $usersSheet)
$ws.Name = "Expression"

# Header row for the new sheet.
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "Expression"
